$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.65
$ws.Cells.Item(2, 9).Value = 5.1
$ws.Cells.Item(2, 16).Value = 2.07
$ws.Cells.Item(2, 17).Value = 1.6
$ws.Cells.Item(2, 18).Value = 5.4
$ws.Cells.Item(2, 19).Value = 6.6
$ws.Cells.Item(2, 20).Value = 8.5
$ws.Cells.Item(2, 21).Value = 11.75
$ws.Cells.Item(2, 22).Value = 15.5
$ws.Cells.Item(2, 24).Value = 7.9
$ws.Cells.Item(2, 25).Value = 6.9
$ws.Cells.Item(2, 26).Value = 21
$ws.Cells.Item(2, 27).Value = 120
$ws.Cells.Item(2, 28).Value = 11.25
$ws.Cells.Item(2, 29).Value = 28
$ws.Cells.Item(2, 30).Value = 17.5
$ws.Cells.Item(2, 31).Value = 100
$ws.Cells.Item(2, 32).Value = 65
$ws.Cells.Item(2, 33).Value = 75
$ws.Cells.Item(5, 7).Value = 4.1
$ws.Cells.Item(5, 9).Value = 1.9
$ws.Cells.Item(6, 12).Value = 1.93
$ws.Cells.Item(6, 13).Value = 1.93
$ws.Cells.Item(8, 14).Value = 1.5
$ws.Cells.Item(8, 15).Value = 2.5
$ws.Cells.Item(8, 18).Value = 5.5
$ws.Cells.Item(8, 26).Value = 21
$ws.Cells.Item(11, 12).Value = 1.85
$ws.Cells.Item(11, 13).Value = 2
$ws.Cells.Item(13, 7).Value = 2.8
$ws.Cells.Item(13, 8).Value = 3.2
$ws.Cells.Item(13, 10).Value = 1.44
$ws.Cells.Item(13, 11).Value = 2.63
$ws.Cells.Item(13, 12).Value = 2.35
$ws.Cells.Item(13, 13).Value = 1.57
$ws.Cells.Item(13, 14).Value = 1.53
$ws.Cells.Item(13, 15).Value = 2.38
$ws.Cells.Item(13, 16).Value = 2.05
$ws.Cells.Item(13, 17).Value = 1.7
$ws.Cells.Item(13, 18).Value = 7
$ws.Cells.Item(13, 22).Value = 26
$ws.Cells.Item(13, 23).Value = 41
$ws.Cells.Item(13, 26).Value = 19
$ws.Cells.Item(13, 27).Value = 67
$ws.Cells.Item(13, 28).Value = 6.5
$ws.Cells.Item(13, 29).Value = 11
$ws.Cells.Item(13, 31).Value = 23
$ws.Cells.Item(13, 33).Value = 41
$ws.Cells.Item(13, 34).Value = 501
$ws.Cells.Item(13, 35).Value = 1.08
$ws.Cells.Item(13, 36).Value = 8
$ws.Cells.Item(16, 7).Value = 4.33
$ws.Cells.Item(16, 8).Value = 4.1
$ws.Cells.Item(16, 9).Value = 1.67
$ws.Cells.Item(16, 12).Value = 1.7
$ws.Cells.Item(16, 13).Value = 2.1
$ws.Cells.Item(16, 16).Value = 1.7
$ws.Cells.Item(16, 17).Value = 2.05
$ws.Cells.Item(17, 7).Value = 3.5
$ws.Cells.Item(17, 10).Value = 1.33
$ws.Cells.Item(17, 11).Value = 3.25
$ws.Cells.Item(17, 18).Value = 9.5
$ws.Cells.Item(17, 21).Value = 34
$ws.Cells.Item(17, 23).Value = 34
$ws.Cells.Item(17, 28).Value = 7.5
$ws.Cells.Item(17, 29).Value = 10
$ws.Cells.Item(17, 31).Value = 21
$ws.Cells.Item(17, 32).Value = 19
$ws.Cells.Item(17, 35).Value = 1.07
$ws.Cells.Item(17, 36).Value = 9
$ws.Cells.Item(19, 7).Value = 1.4
$ws.Cells.Item(19, 8).Value = 4.45
$ws.Cells.Item(19, 9).Value = 5.9
$ws.Cells.Item(19, 12).Value = 1.46
$ws.Cells.Item(19, 13).Value = 2.55
$ws.Cells.Item(19, 16).Value = 1.65
$ws.Cells.Item(19, 17).Value = 2.11
$ws.Cells.Item(19, 18).Value = 8.25
$ws.Cells.Item(19, 19).Value = 7.1
$ws.Cells.Item(19, 20).Value = 7.1
$ws.Cells.Item(19, 21).Value = 8.75
$ws.Cells.Item(19, 22).Value = 8.75
$ws.Cells.Item(19, 23).Value = 16
$ws.Cells.Item(19, 24).Value = 17
$ws.Cells.Item(19, 25).Value = 8.25
$ws.Cells.Item(19, 26).Value = 13
$ws.Cells.Item(19, 27).Value = 40
$ws.Cells.Item(19, 28).Value = 18
$ws.Cells.Item(19, 29).Value = 35
$ws.Cells.Item(19, 30).Value = 15.5
$ws.Cells.Item(19, 31).Value = 90
$ws.Cells.Item(19, 32).Value = 40
$ws.Cells.Item(19, 33).Value = 35
$ws.Cells.Item(19, 34).Value = 200
$ws.Cells.Item(20, 7).Value = 1.11
$ws.Cells.Item(20, 8).Value = 6.8
$ws.Cells.Item(20, 9).Value = 16
$ws.Cells.Item(20, 12).Value = 1.32
$ws.Cells.Item(20, 13).Value = 3.1
$ws.Cells.Item(20, 16).Value = 2.16
$ws.Cells.Item(20, 17).Value = 1.62
$ws.Cells.Item(20, 18).Value = 8.75
$ws.Cells.Item(20, 19).Value = 5.9
$ws.Cells.Item(20, 20).Value = 9.5
$ws.Cells.Item(20, 21).Value = 5.7
$ws.Cells.Item(20, 22).Value = 9
$ws.Cells.Item(20, 23).Value = 26
$ws.Cells.Item(20, 24).Value = 21
$ws.Cells.Item(20, 25).Value = 14
$ws.Cells.Item(20, 26).Value = 27
$ws.Cells.Item(20, 27).Value = 100
$ws.Cells.Item(20, 28).Value = 45
$ws.Cells.Item(20, 29).Value = 150
$ws.Cells.Item(20, 30).Value = 45
$ws.Cells.Item(20, 31).Value = 600
$ws.Cells.Item(20, 32).Value = 200
$ws.Cells.Item(20, 33).Value = 120
$ws.Cells.Item(20, 34).Value = 600
$ws.Cells.Item(21, 7).Value = 5
$ws.Cells.Item(21, 8).Value = 3.9
$ws.Cells.Item(21, 9).Value = 1.55
$ws.Cells.Item(21, 10).Value = 1.17
$ws.Cells.Item(21, 11).Value = 4.5
$ws.Cells.Item(21, 12).Value = 1.55
$ws.Cells.Item(21, 13).Value = 2.29
$ws.Cells.Item(21, 16).Value = 1.63
$ws.Cells.Item(21, 17).Value = 2.15
$ws.Cells.Item(21, 18).Value = 14
$ws.Cells.Item(21, 19).Value = 27
$ws.Cells.Item(21, 20).Value = 13
$ws.Cells.Item(21, 21).Value = 80
$ws.Cells.Item(21, 22).Value = 35
$ws.Cells.Item(21, 23).Value = 35
$ws.Cells.Item(21, 24).Value = 13
$ws.Cells.Item(21, 25).Value = 6.2
$ws.Cells.Item(21, 26).Value = 11
$ws.Cells.Item(21, 27).Value = 45
$ws.Cells.Item(21, 28).Value = 7.4
$ws.Cells.Item(21, 29).Value = 7
$ws.Cells.Item(21, 30).Value = 6.6
$ws.Cells.Item(21, 31).Value = 9.6
$ws.Cells.Item(21, 32).Value = 9.199999999999999
$ws.Cells.Item(21, 33).Value = 17
$ws.Cells.Item(21, 34).Value = 101
$ws.Cells.Item(23, 10).Value = 1.25
$ws.Cells.Item(23, 11).Value = 3.75
$ws.Cells.Item(23, 12).Value = 1.9
$ws.Cells.Item(23, 13).Value = 1.95
$ws.Cells.Item(27, 7).Value = 1.08
$ws.Cells.Item(27, 8).Value = 7.2
$ws.Cells.Item(27, 9).Value = 23
$ws.Cells.Item(27, 12).Value = 1.27
$ws.Cells.Item(27, 13).Value = 3.4
$ws.Cells.Item(27, 16).Value = 2.27
$ws.Cells.Item(27, 17).Value = 1.56
$ws.Cells.Item(27, 18).Value = 9.25
$ws.Cells.Item(27, 19).Value = 6
$ws.Cells.Item(27, 20).Value = 10.25
$ws.Cells.Item(27, 21).Value = 5.5
$ws.Cells.Item(27, 22).Value = 9.25
$ws.Cells.Item(27, 23).Value = 29
$ws.Cells.Item(27, 24).Value = 21
$ws.Cells.Item(27, 25).Value = 16.5
$ws.Cells.Item(27, 26).Value = 32
$ws.Cells.Item(27, 27).Value = 120
$ws.Cells.Item(27, 28).Value = 70
$ws.Cells.Item(27, 29).Value = 300
$ws.Cells.Item(27, 30).Value = 75
$ws.Cells.Item(27, 31).Value = 101
$ws.Cells.Item(27, 32).Value = 400
$ws.Cells.Item(27, 33).Value = 200
$ws.Cells.Item(27, 34).Value = 101
$ws.Cells.Item(30, 12).Value = 2.1
$ws.Cells.Item(30, 13).Value = 1.7
$ws.Cells.Item(30, 14).Value = 1.44
$ws.Cells.Item(30, 15).Value = 2.63
$ws.Cells.Item(31, 10).Value = 1.25
$ws.Cells.Item(31, 11).Value = 3.75
$ws.Cells.Item(31, 12).Value = 1.9
$ws.Cells.Item(31, 13).Value = 1.95
$ws.Cells.Item(31, 14).Value = 1.4
$ws.Cells.Item(32, 7).Value = 1.44
$ws.Cells.Item(32, 14).Value = 1.4
